$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.483.01'
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.066.70'
$ws.Range("E3").Value = '  -0.48%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.91'
$ws.Range("E5").Value = '  +0.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.55'
$ws.Range("E6").Value = '  +1.21%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.541'
$ws.Range("E8").Value = '  -4.54%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.066.97'
$ws.Range("E9").Value = '  -0.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.156'
$ws.Range("E10").Value = '  -0.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.84'
$ws.Range("E11").Value = '  -0.33%  '

$ws.Range("E12").Value = '  -2.76%  '

$ws.Range("E13").Value = '  -1.47%  '

$ws.Range("E14").Value = '  -1.99%  '

$ws.Range("E15").Value = '  +1.38%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.572.26'
$ws.Range("E16").Value = '  -0.55%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.18'
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.415.89'
$ws.Range("E18").Value = '  -0.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.066.93'
$ws.Range("E19").Value = '  -0.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '489.99'
$ws.Range("E20").Value = '  +1.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.42'
$ws.Range("E21").Value = '  -2.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.708'
$ws.Range("E22").Value = '  -1.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.56'
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("E24").Value = '  +4.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.02'
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.90'
$ws.Range("E26").Value = '  -1.69%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.70'
$ws.Range("E27").Value = '  +10.64%  '

$ws.Range("E28").Value = '  +0.05%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.41'
$ws.Range("E29").Value = '  +2.06%  '

$ws.Range("E30").Value = '  +0.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.22'
$ws.Range("E31").Value = '  +1.06%  '

$ws.Range("E32").Value = '  -0.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.44'
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("E34").Value = '  -1.12%  '

$ws.Range("E35").Value = '  +0.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0824'
$ws.Range("E36").Value = '  -3.36%  '

$ws.Range("E37").Value = '  -0.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.00'
$ws.Range("E38").Value = '  -2.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.24'
$ws.Range("E39").Value = '  +0.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.28'
$ws.Range("E40").Value = '  -1.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.62'
$ws.Range("E41").Value = '  +0.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '440.20'
$ws.Range("E42").Value = '  -0.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.293'
$ws.Range("E43").Value = '  +3.07%  '

$ws.Range("E44").Value = '  +2.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0365'
$ws.Range("E45").Value = '  +0.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.850.02'
$ws.Range("E46").Value = '  +1.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '38.91'
$ws.Range("E47").Value = '  -1.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.37'
$ws.Range("E48").Value = '  -0.08%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.34'
$ws.Range("E49").Value = '  +1.05%  '

$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.04%  '

$ws.Range("E51").Value = '  -1.16%  '
